$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student ID in A5 was stored as text; re-enter it as a real number.
$ws.Range("A5").Value = 2590081110

# Remove the stray test/debug values that had been typed into unrelated cells.
$ws.Range("C9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("F21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("E23").ClearContents()

# Reset the view: zoom out to 60% and move the selection to B7.
$excel.ActiveWindow.Zoom = 60
$ws.Range("B7").Select()
